$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.403.98'
$ws.Range('E2').Value = '  +2.72%  '
$ws.Range('D3').Value = '2.503.65'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.89'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.09'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.05%  '
$ws.Range('E7').Value = '  +1.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.55'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0816'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.21%  '
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('D15').Value = '2.894.49'
$ws.Range('E15').Value = '  +2.28%  '
$ws.Range('D16').Value = '2.501.40'
$ws.Range('E16').Value = '  +2.29%  '
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').Value = '47.334.82'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.52%  '
$ws.Range('E20').Value = '  +4.08%  '
$ws.Range('E21').Value = '  +1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +13.65%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '248.57'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.74%  '
$ws.Range('E25').Value = '  +3.90%  '
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('E27').Value = '  -0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.29'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.59%  '
$ws.Range('E29').Value = '  +4.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.36'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.81%  '
$ws.Range('E31').Value = '  +9.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.94'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('E34').Value = '  +1.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0797'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.99%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('E37').Value = '  +5.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.71'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.113'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '121.74'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.37%  '
$ws.Range('E44').Value = '  +2.42%  '
$ws.Range('D45').Value = '2.000.99'
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.59%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  -4.05%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('E50').Value = '  +4.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '56.80'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.16%  '
